# Apply crypto price/volume update for Thu Apr 20 13:29:44 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '29.065.03'
Set-TextValue $ws.Cells.Item(2, 5) '  -1.29%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.987.24'
Set-TextValue $ws.Cells.Item(3, 5) '  -0.41%  '
Set-TextValue $ws.Cells.Item(4, 5) '  +0.32%  '
Set-TextValue $ws.Cells.Item(5, 4) '330.01'
Set-TextValue $ws.Cells.Item(5, 5) '  +0.14%  '
Set-TextValue $ws.Cells.Item(6, 4) '1.010'
Set-TextValue $ws.Cells.Item(6, 5) '  +0.24%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.4976'
Set-TextValue $ws.Cells.Item(7, 5) '  -0.65%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.4195'
Set-TextValue $ws.Cells.Item(8, 5) '  -0.78%  '
Set-TextValue $ws.Cells.Item(9, 5) '  +1.93%  '
Set-TextValue $ws.Cells.Item(10, 4) '0.09318'
Set-TextValue $ws.Cells.Item(10, 5) '  +4.33%  '
Set-TextValue $ws.Cells.Item(11, 4) '1.094'
Set-TextValue $ws.Cells.Item(11, 5) '  -2.43%  '
Set-TextValue $ws.Cells.Item(12, 4) '23.15'
Set-TextValue $ws.Cells.Item(12, 5) '  -0.41%  '
Set-TextValue $ws.Cells.Item(13, 4) '2.014.67'
Set-TextValue $ws.Cells.Item(13, 5) '  +4.51%  '
Set-TextValue $ws.Cells.Item(14, 4) '7.977'
Set-TextValue $ws.Cells.Item(14, 5) '  -1.36%  '
Set-TextValue $ws.Cells.Item(15, 4) '6.442'
Set-TextValue $ws.Cells.Item(15, 5) '  -1.23%  '
Set-TextValue $ws.Cells.Item(16, 5) '  +0.38%  '
$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Cells.Item(17, 4) '0.00001112'
Set-TextValue $ws.Cells.Item(17, 5) '  +0.26%  '
$ws.Cells.Item(18, 2).Value = 'Litecoin'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Cells.Item(18, 4) '92.40'
Set-TextValue $ws.Cells.Item(18, 5) '  -3.53%  '
Set-TextValue $ws.Cells.Item(19, 4) '0.06751'
Set-TextValue $ws.Cells.Item(19, 5) '  +2.00%  '
Set-TextValue $ws.Cells.Item(20, 4) '19.55'
Set-TextValue $ws.Cells.Item(20, 5) '  -0.58%  '
Set-TextValue $ws.Cells.Item(21, 4) '1.010'
Set-TextValue $ws.Cells.Item(21, 5) '  +0.35%  '
Set-TextValue $ws.Cells.Item(22, 4) '5.972'
Set-TextValue $ws.Cells.Item(22, 5) '  +0.02%  '
Set-TextValue $ws.Cells.Item(23, 4) '29.064.96'
Set-TextValue $ws.Cells.Item(23, 5) '  -1.36%  '
Set-TextValue $ws.Cells.Item(24, 4) '11.98'
Set-TextValue $ws.Cells.Item(24, 5) '  +0.76%  '
Set-TextValue $ws.Cells.Item(25, 4) '2.286'
Set-TextValue $ws.Cells.Item(25, 5) '  +1.40%  '
Set-TextValue $ws.Cells.Item(26, 4) '2.245.88'
Set-TextValue $ws.Cells.Item(26, 5) '  +3.17%  '
Set-TextValue $ws.Cells.Item(27, 5) '  +0.71%  '
Set-TextValue $ws.Cells.Item(28, 4) '156.55'
Set-TextValue $ws.Cells.Item(28, 5) '  -1.07%  '
Set-TextValue $ws.Cells.Item(29, 4) '6.269'
Set-TextValue $ws.Cells.Item(29, 5) '  -4.38%  '
Set-TextValue $ws.Cells.Item(30, 4) '2.260'
Set-TextValue $ws.Cells.Item(30, 5) '  -2.92%  '
Set-TextValue $ws.Cells.Item(31, 4) '127.46'
Set-TextValue $ws.Cells.Item(31, 5) '  -0.28%  '
Set-TextValue $ws.Cells.Item(32, 4) '1.047'
Set-TextValue $ws.Cells.Item(32, 5) '  -0.07%  '
Set-TextValue $ws.Cells.Item(33, 4) '0.09840'
Set-TextValue $ws.Cells.Item(33, 5) '  -1.05%  '
Set-TextValue $ws.Cells.Item(34, 4) '1.532'
Set-TextValue $ws.Cells.Item(34, 5) '  -2.35%  '
Set-TextValue $ws.Cells.Item(35, 4) '5.811'
Set-TextValue $ws.Cells.Item(35, 5) '  -0.62%  '
Set-TextValue $ws.Cells.Item(36, 4) '3.746'
Set-TextValue $ws.Cells.Item(36, 5) '  -1.01%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.02424'
Set-TextValue $ws.Cells.Item(37, 5) '  -1.56%  '
Set-TextValue $ws.Cells.Item(38, 4) '1.316'
Set-TextValue $ws.Cells.Item(38, 5) '  +2.28%  '
Set-TextValue $ws.Cells.Item(39, 4) '9.058'
Set-TextValue $ws.Cells.Item(39, 5) '  -5.37%  '
Set-TextValue $ws.Cells.Item(40, 5) '  +1.02%  '
Set-TextValue $ws.Cells.Item(41, 4) '0.6481'
Set-TextValue $ws.Cells.Item(41, 5) '  -0.63%  '
Set-TextValue $ws.Cells.Item(42, 5) '  -1.74%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.1992'
Set-TextValue $ws.Cells.Item(43, 5) '  -3.58%  '
Set-TextValue $ws.Cells.Item(44, 5) '  +0.25%  '
Set-TextValue $ws.Cells.Item(45, 4) '0.6220'
Set-TextValue $ws.Cells.Item(45, 5) '  -1.73%  '
Set-TextValue $ws.Cells.Item(46, 4) '1.347'
Set-TextValue $ws.Cells.Item(46, 5) '  +6.33%  '
Set-TextValue $ws.Cells.Item(47, 4) '13.44'
Set-TextValue $ws.Cells.Item(47, 5) '  +0.44%  '
Set-TextValue $ws.Cells.Item(48, 4) '2.183'
Set-TextValue $ws.Cells.Item(48, 5) '  -1.00%  '
Set-TextValue $ws.Cells.Item(49, 4) '3.489'
Set-TextValue $ws.Cells.Item(49, 5) '  -1.10%  '
Set-TextValue $ws.Cells.Item(50, 4) '0.00000000340'
Set-TextValue $ws.Cells.Item(50, 5) '  +3.11%  '
Set-TextValue $ws.Cells.Item(51, 4) '0.06985'
Set-TextValue $ws.Cells.Item(51, 5) '  -0.15%  '
